# Refresh the CodeSystem metadata sheet to reflect a new publish run:
#   Status        draft -> active
#   Experimental  (blank) -> false
#   Date          2025-05-21T14:22:51+00:00 -> 2025-06-13T15:45:04+00:00
#   Case Sensitive (blank) -> true

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status (row 6)
$ws.Cells.Item(6, 2).Value = "active"

# Date (row 8)
$ws.Cells.Item(8, 2).Value = "2025-06-13T15:45:04+00:00"

# Experimental (row 7): must land as literal text "false", not the
# Boolean FALSE that a bare Value assignment would auto-coerce it to.
# Typing it with a leading apostrophe forces text entry (as it would in
# the Excel UI); we then restore the plain data-row formatting via
# PasteSpecial(Formats) from a neighboring cell so the quote-prefix
# formatting doesn't stick to this cell.
$expCell = $ws.Cells.Item(7, 2)
$expCell.Value = "'false"
$ws.Cells.Item(8, 2).Copy()
$expCell.PasteSpecial($xlPasteFormats)

# Case Sensitive (row 15): same literal-text trick for "true".
$csCell = $ws.Cells.Item(15, 2)
$csCell.Value = "'true"
$ws.Cells.Item(14, 2).Copy()
$csCell.PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
